$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "abcdef"
$ws.Range("B3").Value = "password"

$ws.Range("B3").Select()
